$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F4").Value = 0
$ws.Range("F6").Value = 4
$ws.Range("F13").Value = -2
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = 4
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 2
$ws.Range("F22").Value = 5
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("F25").Value = -2
$ws.Range("F26").Value = -5
$ws.Range("F27").Value = 1
$ws.Range("F28").Value = -3
$ws.Range("F29").Value = -1
$ws.Range("F30").Value = -2
$ws.Range("F31").Value = 2
$ws.Range("F32").Value = 2
$ws.Range("F33").Value = 2
$ws.Range("F34").Value = -1
$ws.Range("F35").Value = -2
$ws.Range("F37").Value = -1
$ws.Range("F42").Value = 6
$ws.Range("F43").Value = -1
$ws.Range("F44").Value = -4
$ws.Range("F45").Value = -2
$ws.Range("F46").Value = 5
$ws.Range("F48").Value = 9
$ws.Range("F49").Value = -5
